# Inserts a new weekly price record (Macroferia Regional de Talca - Cilantro)
# as row 62, pushing the two previously-last rows down to 63/64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62 (shifts old rows 62 and 63 down to 63/64).
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record's data.
$ws.Cells.Item(62, 1).Value  = 5
$ws.Cells.Item(62, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(62, 3).Value  = "Maule"
$ws.Cells.Item(62, 4).Value  = 44826
$ws.Cells.Item(62, 5).Value  = 7
$ws.Cells.Item(62, 6).Value  = 100112040
$ws.Cells.Item(62, 7).Value  = "Cilantro"
$ws.Cells.Item(62, 8).Value  = "Sin especificar"
$ws.Cells.Item(62, 9).Value  = "Primera"
$ws.Cells.Item(62, 10).Value = 150
$ws.Cells.Item(62, 11).Value = 8000
$ws.Cells.Item(62, 12).Value = 8000
$ws.Cells.Item(62, 13).Value = 8000
$ws.Cells.Item(62, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(62, 15).Value = "Región del Maule"
$ws.Cells.Item(62, 16).Value = 222
$ws.Cells.Item(62, 17).Value = 36
$ws.Cells.Item(62, 18).Value = "Hortaliza"
